$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "276.03"
Set-TextValue $ws.Range("G2") "15"

# Row 3
Set-TextValue $ws.Range("D3") "20.90"
Set-TextValue $ws.Range("G3") "15"

# Row 4
Set-TextValue $ws.Range("G4") "15"

# Row 5
Set-TextValue $ws.Range("D5") "0.06199"
Set-TextValue $ws.Range("G5") "15"

# Row 6
Set-TextValue $ws.Range("D6") "3.578"
Set-TextValue $ws.Range("G6") "15"

# Row 7
Set-TextValue $ws.Range("G7") "15"

# Row 8
Set-TextValue $ws.Range("D8") "1.480"
Set-TextValue $ws.Range("G8") "15"

# Row 9
Set-TextValue $ws.Range("D9") "0.8208"
Set-TextValue $ws.Range("G9") "15"

# Row 10
Set-TextValue $ws.Range("D10") "0.1637"
Set-TextValue $ws.Range("G10") "15"

# Row 11
Set-TextValue $ws.Range("D11") "0.08292"
Set-TextValue $ws.Range("G11") "15"

# Row 12
Set-TextValue $ws.Range("D12") "0.03491"
Set-TextValue $ws.Range("G12") "15"

# Row 13
Set-TextValue $ws.Range("D13") "0.03126"
Set-TextValue $ws.Range("G13") "15"

# Row 14
Set-TextValue $ws.Range("D14") "0.09129"
Set-TextValue $ws.Range("G14") "15"

# Row 15
Set-TextValue $ws.Range("D15") "3.773"
Set-TextValue $ws.Range("G15") "15"

# Row 16
Set-TextValue $ws.Range("D16") "0.001625"
Set-TextValue $ws.Range("G16") "15"

# Row 17
Set-TextValue $ws.Range("G17") "15"

# Row 18
Set-TextValue $ws.Range("D18") "0.006396"
Set-TextValue $ws.Range("G18") "15"

# Row 19
Set-TextValue $ws.Range("G19") "15"

# Row 20
Set-TextValue $ws.Range("D20") "0.001067"
Set-TextValue $ws.Range("G20") "15"

# Row 21
Set-TextValue $ws.Range("G21") "15"

# Row 22
Set-TextValue $ws.Range("D22") "3.804"
Set-TextValue $ws.Range("G22") "15"

# Row 23
Set-TextValue $ws.Range("D23") "2.321"
Set-TextValue $ws.Range("G23") "15"

# Row 24
Set-TextValue $ws.Range("D24") "0.01384"
Set-TextValue $ws.Range("G24") "15"

# Row 25
Set-TextValue $ws.Range("D25") "0.3386"
Set-TextValue $ws.Range("G25") "15"

# Row 26
Set-TextValue $ws.Range("D26") "0.1201"
Set-TextValue $ws.Range("G26") "15"

# Row 27
Set-TextValue $ws.Range("G27") "15"

# Row 28
Set-TextValue $ws.Range("G28") "15"

# Row 29
Set-TextValue $ws.Range("G29") "15"

# Row 30
Set-TextValue $ws.Range("G30") "15"

# Row 31
Set-TextValue $ws.Range("G31") "15"

# Row 32
Set-TextValue $ws.Range("G32") "15"

# Row 33
Set-TextValue $ws.Range("G33") "15"

# Row 34
Set-TextValue $ws.Range("G34") "15"

# Row 35
Set-TextValue $ws.Range("G35") "15"

# Row 36
Set-TextValue $ws.Range("G36") "15"

# Row 37
Set-TextValue $ws.Range("G37") "15"

# Row 38
Set-TextValue $ws.Range("G38") "15"

# Row 39
Set-TextValue $ws.Range("G39") "15"

# Row 40
Set-TextValue $ws.Range("D40") "0.04668"
Set-TextValue $ws.Range("G40") "15"

# Row 41
Set-TextValue $ws.Range("D41") "0.007036"
Set-TextValue $ws.Range("G41") "15"

# Row 42
Set-TextValue $ws.Range("B42") "BKEXToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1105"
Set-TextValue $ws.Range("E42") "41BKEXTokenBKK"
Set-TextValue $ws.Range("G42") "15"

# Row 43
Set-TextValue $ws.Range("B43") "CEJI"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003523"
Set-TextValue $ws.Range("E43") "42CEJICEJI"
Set-TextValue $ws.Range("G43") "15"

# Row 44
Set-TextValue $ws.Range("D44") "0.01116"
Set-TextValue $ws.Range("G44") "15"

# Row 45
Set-TextValue $ws.Range("D45") "0.00006263"
Set-TextValue $ws.Range("G45") "15"

# Row 46
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("G46") "15"

# Row 47
Set-TextValue $ws.Range("D47") "0.7921"
Set-TextValue $ws.Range("G47") "15"

# Row 48
Set-TextValue $ws.Range("D48") "0.002304"
Set-TextValue $ws.Range("G48") "15"

# Row 49
Set-TextValue $ws.Range("G49") "15"

# Row 50
Set-TextValue $ws.Range("G50") "15"

# Row 51
Set-TextValue $ws.Range("G51") "15"
